$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New series label (row 1, shared string "210 Packets") ---
$ws.Range("V1").Value = "210 Packets"

# --- New sub-header row (row 2): Acc / Loss / Time ---
$ws.Range("U2").Value = "Acc"
$ws.Range("V2").Value = "Loss"
$ws.Range("W2").Value = "Time"

# --- New data column (210 Packets run): Q=Acc, R=Loss, S=Time, rows 3-51 ---
$ws.Range("Q3").Value = 89.072024822235093
$ws.Range("R3").Value = 0.26912349552781201
$ws.Range("S3").Value = 76.8498628139495
$ws.Range("Q4").Value = 88.799297809600802
$ws.Range("R4").Value = 0.25287529295787198
$ws.Range("S4").Value = 68.147766590118394
$ws.Range("Q5").Value = 89.067280292510901
$ws.Range("R5").Value = 0.25610320560331401
$ws.Range("S5").Value = 67.9189066886901
$ws.Range("Q6").Value = 88.770842552185002
$ws.Range("R6").Value = 0.263082958474816
$ws.Range("S6").Value = 68.626559019088702
$ws.Range("Q7").Value = 88.4909987449646
$ws.Range("R7").Value = 0.25077301616098102
$ws.Range("S7").Value = 67.212399482726994
$ws.Range("Q8").Value = 89.019846916198702
$ws.Range("R8").Value = 0.25939821293550802
$ws.Range("S8").Value = 68.903663635253906
$ws.Range("Q9").Value = 88.588231801986694
$ws.Range("R9").Value = 0.25101132475818999
$ws.Range("S9").Value = 67.710575819015503
$ws.Range("Q10").Value = 89.026963710784898
$ws.Range("R10").Value = 0.25330214329852102
$ws.Range("S10").Value = 68.029569864273
$ws.Range("Q11").Value = 89.674389362335205
$ws.Range("R11").Value = 0.246350247479621
$ws.Range("S11").Value = 69.242267370223999
$ws.Range("Q12").Value = 89.349490404129
$ws.Range("R12").Value = 0.25696080555966999
$ws.Range("S12").Value = 68.885797739028902
$ws.Range("Q13").Value = 89.223802089691105
$ws.Range("R13").Value = 0.25985198954120797
$ws.Range("S13").Value = 68.990481376647907
$ws.Range("Q14").Value = 89.278346300125094
$ws.Range("R14").Value = 0.25125204967355402
$ws.Range("S14").Value = 69.237408876418996
$ws.Range("Q15").Value = 89.294946193694997
$ws.Range("R15").Value = 0.25244044060549198
$ws.Range("S15").Value = 67.954193115234304
$ws.Range("Q16").Value = 89.913916587829505
$ws.Range("R16").Value = 0.240826629251888
$ws.Range("S16").Value = 68.948252916336003
$ws.Range("Q17").Value = 89.098107814788804
$ws.Range("R17").Value = 0.255289133261347
$ws.Range("S17").Value = 68.126592636108398
$ws.Range("Q18").Value = 89.522612094879094
$ws.Range("R18").Value = 0.246590712003668
$ws.Range("S18").Value = 68.684085369110093
$ws.Range("Q19").Value = 88.875186443328801
$ws.Range("R19").Value = 0.26004882433884602
$ws.Range("S19").Value = 68.964230775833101
$ws.Range("Q20").Value = 89.219057559966998
$ws.Range("R20").Value = 0.25694732281784199
$ws.Range("S20").Value = 68.421115875244098
$ws.Range("Q21").Value = 89.534467458724905
$ws.Range("R21").Value = 0.25473587350469601
$ws.Range("S21").Value = 69.155242919921804
$ws.Range("Q22").Value = 89.394551515579195
$ws.Range("R22").Value = 0.25280641418191602
$ws.Range("S22").Value = 68.854432821273804
$ws.Range("Q23").Value = 88.851469755172701
$ws.Range("R23").Value = 0.26166377226755
$ws.Range("S23").Value = 68.335008144378605
$ws.Range("Q24").Value = 89.278346300125094
$ws.Range("R24").Value = 0.25392871053295302
$ws.Range("S24").Value = 68.157465219497595
$ws.Range("Q25").Value = 88.450682163238497
$ws.Range("R25").Value = 0.26895801092282801
$ws.Range("S25").Value = 67.761798143386798
$ws.Range("Q26").Value = 89.733678102493201
$ws.Range("R26").Value = 0.248675107689834
$ws.Range("S26").Value = 68.530229330062795
$ws.Range("Q27").Value = 89.368462562561007
$ws.Range("R27").Value = 0.25674521599223599
$ws.Range("S27").Value = 68.778434991836505
$ws.Range("Q28").Value = 89.038819074630695
$ws.Range("R28").Value = 0.25517505047163902
$ws.Range("S28").Value = 67.442076683044405
$ws.Range("Q29").Value = 89.510756731033297
$ws.Range("R29").Value = 0.24333109311892101
$ws.Range("S29").Value = 68.255703687667804
$ws.Range("Q30").Value = 89.543956518173204
$ws.Range("R30").Value = 0.245791154289593
$ws.Range("S30").Value = 69.228202104568396
$ws.Range("Q31").Value = 88.965302705764699
$ws.Range("R31").Value = 0.26215698969425799
$ws.Range("S31").Value = 68.481005430221501
$ws.Range("Q32").Value = 89.152652025222693
$ws.Range("R32").Value = 0.251220728692363
$ws.Range("S32").Value = 67.891291379928504
$ws.Range("Q33").Value = 89.098107814788804
$ws.Range("R33").Value = 0.24443516779753599
$ws.Range("S33").Value = 68.9181227684021
$ws.Range("Q34").Value = 88.851469755172701
$ws.Range("R34").Value = 0.26175096853258201
$ws.Range("S34").Value = 68.831122875213595
$ws.Range("Q35").Value = 89.143168926238999
$ws.Range("R35").Value = 0.250213219927652
$ws.Range("S35").Value = 69.915167093276906
$ws.Range("Q36").Value = 89.747905731201101
$ws.Range("R36").Value = 0.251498217096777
$ws.Range("S36").Value = 69.450585126876803
$ws.Range("Q37").Value = 88.804042339324894
$ws.Range("R37").Value = 0.26033363638799401
$ws.Range("S37").Value = 69.064718723297105
$ws.Range("Q38").Value = 88.493371009826603
$ws.Range("R38").Value = 0.25170568995396603
$ws.Range("S38").Value = 68.731973648071204
$ws.Range("Q39").Value = 88.372421264648395
$ws.Range("R39").Value = 0.25919652122440601
$ws.Range("S39").Value = 68.462583303451495
$ws.Range("Q40").Value = 89.112341403961096
$ws.Range("R40").Value = 0.252560794765531
$ws.Range("S40").Value = 68.509573459625202
$ws.Range("Q41").Value = 89.067280292510901
$ws.Range("R41").Value = 0.25214308744849601
$ws.Range("S41").Value = 69.308448314666705
$ws.Range("Q42").Value = 89.944744110107393
$ws.Range("R42").Value = 0.24519696594476201
$ws.Range("S42").Value = 69.181425809860201
$ws.Range("Q43").Value = 88.661748170852604
$ws.Range("R43").Value = 0.25244328693842799
$ws.Range("S43").Value = 67.910869598388601
$ws.Range("Q44").Value = 89.019846916198702
$ws.Range("R44").Value = 0.25686037386819499
$ws.Range("S44").Value = 67.286752700805593
$ws.Range("Q45").Value = 88.896530866622896
$ws.Range("R45").Value = 0.26274500668609302
$ws.Range("S45").Value = 68.270944833755493
$ws.Range("Q46").Value = 89.019846916198702
$ws.Range("R46").Value = 0.25729004903374297
$ws.Range("S46").Value = 69.058609247207599
$ws.Range("Q47").Value = 88.611948490142794
$ws.Range("R47").Value = 0.26704773569128498
$ws.Range("S47").Value = 67.666248798370304
$ws.Range("Q48").Value = 89.968460798263493
$ws.Range("R48").Value = 0.246782750153139
$ws.Range("S48").Value = 68.321039438247595
$ws.Range("Q49").Value = 89.306801557540894
$ws.Range("R49").Value = 0.25290859650150699
$ws.Range("S49").Value = 68.892151355743394
$ws.Range("Q50").Value = 88.792181015014606
$ws.Range("R50").Value = 0.25707833496460297
$ws.Range("S50").Value = 68.697601079940796
$ws.Range("Q51").Value = 89.145541191101003
$ws.Range("R51").Value = 0.25451644660448802
$ws.Range("S51").Value = 68.886629104614201

# --- Move the selection to match where the author left off editing ---
$ws.Range("U6").Select()
